# Updates the "Improved Wind Production Forecast" data: shifts the timestamp
# series forward by one week (45751 -> 45758) and refreshes the forecasted
# consumption values (column A) for rows 2-97 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Forecasted Consumption (MW)" values for A2:A97
$consumption = @(
    6070, 6020, 5970, 5940, 5910, 5880, 5860, 5850, 5850, 5860, 5870, 5880,
    5890, 5910, 5930, 5960, 6010, 6070, 6160, 6270, 6400, 6540, 6690, 6850,
    7000, 7150, 7270, 7370, 7440, 7480, 7500, 7480, 7440, 7380, 7310, 7220,
    7140, 7060, 6980, 6910, 6850, 6790, 6740, 6690, 6640, 6590, 6550, 6520,
    6490, 6470, 6460, 6450, 6460, 6460, 6470, 6470, 6460, 6460, 6450, 6450,
    6450, 6480, 6510, 6570, 6630, 6700, 6770, 6840, 6920, 7000, 7100, 7200,
    7330, 7450, 7570, 7680, 7760, 7830, 7860, 7880, 7840, 7780, 7660, 7540,
    7380, 7220, 7050, 6880, 6720, 6570, 6420, 6280, 6150, 6090, 6040, 5980
)

# New "Timestamp" values (date serials) for B2:B97
$timestamps = @(
    45758, 45758.01041666666, 45758.02083333334, 45758.03125, 45758.04166666666, 45758.05208333334, 45758.0625, 45758.07291666666, 45758.08333333334, 45758.09375, 45758.10416666666, 45758.11458333334,
    45758.125, 45758.13541666666, 45758.14583333334, 45758.15625, 45758.16666666666, 45758.17708333334, 45758.1875, 45758.19791666666, 45758.20833333334, 45758.21875, 45758.22916666666, 45758.23958333334,
    45758.25, 45758.26041666666, 45758.27083333334, 45758.28125, 45758.29166666666, 45758.30208333334, 45758.3125, 45758.32291666666, 45758.33333333334, 45758.34375, 45758.35416666666, 45758.36458333334,
    45758.375, 45758.38541666666, 45758.39583333334, 45758.40625, 45758.41666666666, 45758.42708333334, 45758.4375, 45758.44791666666, 45758.45833333334, 45758.46875, 45758.47916666666, 45758.48958333334,
    45758.5, 45758.51041666666, 45758.52083333334, 45758.53125, 45758.54166666666, 45758.55208333334, 45758.5625, 45758.57291666666, 45758.58333333334, 45758.59375, 45758.60416666666, 45758.61458333334,
    45758.625, 45758.63541666666, 45758.64583333334, 45758.65625, 45758.66666666666, 45758.67708333334, 45758.6875, 45758.69791666666, 45758.70833333334, 45758.71875, 45758.72916666666, 45758.73958333334,
    45758.75, 45758.76041666666, 45758.77083333334, 45758.78125, 45758.79166666666, 45758.80208333334, 45758.8125, 45758.82291666666, 45758.83333333334, 45758.84375, 45758.85416666666, 45758.86458333334,
    45758.875, 45758.88541666666, 45758.89583333334, 45758.90625, 45758.91666666666, 45758.92708333334, 45758.9375, 45758.94791666666, 45758.95833333334, 45758.96875, 45758.97916666666, 45758.98958333334
)

for ($i = 0; $i -lt $consumption.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $consumption[$i]
    $ws.Cells.Item($row, 2).Value = $timestamps[$i]
}
